# Apply the "synced files as of 10.27.24" update to the Contribution Log.
# This adds a new contribution-log entry (row 8) describing the Raspberry Pi
# remote configuration work, and updates the active selection to reflect
# where the author left the cursor after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entry (row 8) -------------------------------------------------
# Team member (dropdown-validated list value)
$ws.Range("A8").Value = "Erik Meurrens"

# Item name / description (new shared strings get appended automatically)
$ws.Range("B8").Value = "Raspberry Pi remote configuration"
$ws.Range("C8").Value = "Configured SSH/VNC and Jupyter Notebook to be able to access board, using the UF VPN, while it is connected to the UF WiFi network"

# Start / End dates -- copy the existing date cell's number formatting first
# so the new cells reuse the same style (m/d/yyyy) instead of creating a new
# style definition, then set the actual serial date values.
$ws.Range("D5").Copy()
$ws.Range("D8:E8").PasteSpecial(-4122)
$ws.Range("D8").Value = 45580
$ws.Range("E8").Value = 45585

# --- Row heights ------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(3).RowHeight = 44.4
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 28.8

# --- Selection ---------------------------------------------------------------
$ws.Range("E9").Select()
